$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 25 de Septiembre de 2020 a las 20:36"
$ws.Range("B4").Value = 7206769
$ws.Range("C4").Value = 21298
$ws.Range("D4").Value = 4459624
$ws.Range("E4").Value = 2539222
$ws.Range("G4").Value = 385
$ws.Range("H4").Value = 207923
$ws.Range("B5").Value = 5897227
$ws.Range("C5").Value = 81124
$ws.Range("D5").Value = 4836208
$ws.Range("E5").Value = 967645
$ws.Range("G5").Value = 1057
$ws.Range("H5").Value = 93374
$ws.Range("B6").Value = 4667384
$ws.Range("C6").Value = 7475
$ws.Range("E6").Value = 503555
$ws.Range("G6").Value = 157
$ws.Range("H6").Value = 140040
$ws.Range("A10").Value = "España"
$ws.Range("B10").Value = 735198
$ws.Range("C10").Value = 4122
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("G10").Value = 114
$ws.Range("H10").Value = 31232
$ws.Range("A11").Value = "Mexico"
$ws.Range("B11").Value = 715457
$ws.Range("C11").Value = 5408
$ws.Range("D11").Value = 514760
$ws.Range("E11").Value = 125258
$ws.Range("G11").Value = 490
$ws.Range("H11").Value = 75439
$ws.Range("B14").Value = 513034
$ws.Range("C14").Value = 15797
$ws.Range("D14").Value = 94891
$ws.Range("E14").Value = 386482
$ws.Range("G14").Value = 150
$ws.Range("H14").Value = 31661
$ws.Range("B17").Value = 423236
$ws.Range("C17").Value = 6874
$ws.Range("G17").Value = 34
$ws.Range("H17").Value = 41936
$ws.Range("B25").Value = 283444
$ws.Range("C25").Value = 2099
$ws.Range("E25").Value = 24416
$ws.Range("G25").Value = 9
$ws.Range("H25").Value = 9528
$ws.Range("B62").Value = 50754
$ws.Range("C62").Value = 175
$ws.Range("D62").Value = 35654
$ws.Range("E62").Value = 13393
$ws.Range("G62").Value = 4
$ws.Range("H62").Value = 1707
$ws.Range("B73").Value = 34315
$ws.Range("C73").Value = 320
$ws.Range("E73").Value = 9154
$ws.Range("B96").Value = 10835
$ws.Range("C96").Value = 95
$ws.Range("D96").Value = 8569
$ws.Range("E96").Value = 2147
$ws.Range("A111").Value = "Mozambique"
$ws.Range("B111").Value = 7589
$ws.Range("C111").Value = 190
$ws.Range("D111").Value = 4649
$ws.Range("E111").Value = 2887
$ws.Range("G111").Value = 2
$ws.Range("H111").Value = 53
$ws.Range("A112").Value = "Mauritania"
$ws.Range("B112").Value = 7433
$ws.Range("D112").Value = 7052
$ws.Range("E112").Value = 220
$ws.Range("H112").Value = 161
$ws.Range("B115").Value = 5764
$ws.Range("C115").Value = 17
$ws.Range("D115").Value = 4178
$ws.Range("E115").Value = 1407
$ws.Range("B119").Value = 5399
$ws.Range("C119").Value = 24
$ws.Range("D119").Value = 4767
$ws.Range("E119").Value = 524
$ws.Range("A137").Value = "Reunion"
$ws.Range("B137").Value = 3685
$ws.Range("C137").Value = 184
$ws.Range("D137").Value = 2819
$ws.Range("E137").Value = 855
$ws.Range("H137").Value = 11
$ws.Range("A138").Value = "Gambia"
$ws.Range("B138").Value = 3555
$ws.Range("C138").Value = 3
$ws.Range("D138").Value = 2034
$ws.Range("E138").Value = 1411
$ws.Range("H138").Value = 110
$ws.Range("A139").Value = "Mayotte"
$ws.Range("B139").Value = 3541
$ws.Range("C139").Value = 0
$ws.Range("D139").Value = 2964
$ws.Range("E139").Value = 537
$ws.Range("H139").Value = 40
$ws.Range("A140").Value = "Tailandia"
$ws.Range("B140").Value = 3519
$ws.Range("C140").Value = 3
$ws.Range("D140").Value = 3360
$ws.Range("E140").Value = 100
$ws.Range("H140").Value = 59
$ws.Range("B142").Value = 3345
$ws.Range("C142").Value = 12
$ws.Range("E142").Value = 174
$ws.Range("B185").Value = 329
$ws.Range("C185").Value = 14
$ws.Range("D185").Value = 124
$ws.Range("E185").Value = 204
$ws.Range("A215").Value = "Montserrat"
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1
$ws.Range("A216").Value = "Islas Malvinas"
$ws.Range("D216").Value = 13
$ws.Range("H216").Value = 0
